$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JSS 3D")

$ws.Range("A8").Value = "2026-02-10 19:06:54"
$ws.Range("B8").Value = "Mahmud Alhaji Hassan"
$ws.Range("C8").Value = "Number 5"
$ws.Range("D8").Value = 9
